$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new Scrum meeting attendance row (row 13)
$ws.Range("B13").Value = "9/15 /4:15"
$ws.Range("C13").Value = "Google Hangout"
$ws.Range("D13").Value = "A"
$ws.Range("E13").Value = "A"
$ws.Range("F13").Value = "A"
$ws.Range("G13").Value = "A"
$ws.Range("H13").Value = "A"
$ws.Range("I13").Value = "A"

# C13 picks up the same "no top border / bold" look used by the
# meeting-place cells directly above it (C6:C12)
$ws.Range("C13").Font.Bold = $true
$ws.Range("C13").Borders(8).LineStyle = -4142

# Move the active selection to the newly filled-in row
$ws.Range("B13").Select()
